# Auto-generated edit script applying numeric corrections to the
# Seraph_Profits leve-profit tables across all 8 job sheets,
# per the scheduled price-refresh run (chore: update Sheets via scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 4207.9165
$ws.Cells.Item(17, 10).Value = 4207.9165
$ws.Cells.Item(17, 12).Value = 12623.7495
$ws.Cells.Item(17, 14).Value = -12959.7495
$ws.Cells.Item(86, 8).Value = 4058.8
$ws.Cells.Item(86, 9).Value = 2950
$ws.Cells.Item(86, 10).Value = 4798
$ws.Cells.Item(86, 11).Value = 2950
$ws.Cells.Item(86, 12).Value = 4798
$ws.Cells.Item(86, 13).Value = -1827
$ws.Cells.Item(86, 14).Value = -7044
$ws.Cells.Item(89, 8).Value = 4058.8
$ws.Cells.Item(89, 9).Value = 2950
$ws.Cells.Item(89, 10).Value = 4798
$ws.Cells.Item(89, 11).Value = 14750
$ws.Cells.Item(89, 12).Value = 23990
$ws.Cells.Item(89, 13).Value = -9134
$ws.Cells.Item(89, 14).Value = -35222
$ws.Cells.Item(132, 8).Value = 2261.125
$ws.Cells.Item(132, 9).Value = 2033.4
$ws.Cells.Item(132, 11).Value = 6100.200000000001
$ws.Cells.Item(132, 13).Value = -3570.200000000001
$ws.Cells.Item(137, 8).Value = 1838.591
$ws.Cells.Item(137, 9).Value = 1721.2354
$ws.Cells.Item(137, 11).Value = 5163.706200000001
$ws.Cells.Item(137, 13).Value = -2613.706200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 419465.78
$ws.Cells.Item(122, 9).Value = 669423
$ws.Cells.Item(122, 10).Value = 2870.4443
$ws.Cells.Item(122, 11).Value = 2008269
$ws.Cells.Item(122, 12).Value = 8611.332900000001
$ws.Cells.Item(122, 13).Value = -2005819
$ws.Cells.Item(122, 14).Value = -13511.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 47355.184
$ws.Cells.Item(99, 9).Value = 68573.53
$ws.Cells.Item(99, 11).Value = 68573.53
$ws.Cells.Item(99, 13).Value = -67075.53
$ws.Cells.Item(134, 8).Value = 1921.7273
$ws.Cells.Item(134, 9).Value = 1645.4
$ws.Cells.Item(134, 10).Value = 2785.25
$ws.Cells.Item(134, 11).Value = 4936.200000000001
$ws.Cells.Item(134, 12).Value = 8355.75
$ws.Cells.Item(134, 13).Value = -2401.200000000001
$ws.Cells.Item(134, 14).Value = -13425.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3690.0444
$ws.Cells.Item(31, 9).Value = 3073.4736
$ws.Cells.Item(31, 11).Value = 3073.4736
$ws.Cells.Item(31, 13).Value = -2778.4736
$ws.Cells.Item(34, 8).Value = 3690.0444
$ws.Cells.Item(34, 9).Value = 3073.4736
$ws.Cells.Item(34, 11).Value = 3073.4736
$ws.Cells.Item(34, 13).Value = -2871.4736
$ws.Cells.Item(86, 8).Value = 9762.272000000001
$ws.Cells.Item(86, 9).Value = 8599.200000000001
$ws.Cells.Item(86, 10).Value = 10731.5
$ws.Cells.Item(86, 11).Value = 8599.200000000001
$ws.Cells.Item(86, 12).Value = 10731.5
$ws.Cells.Item(86, 13).Value = -7476.200000000001
$ws.Cells.Item(86, 14).Value = -12977.5
$ws.Cells.Item(89, 8).Value = 9762.272000000001
$ws.Cells.Item(89, 9).Value = 8599.200000000001
$ws.Cells.Item(89, 10).Value = 10731.5
$ws.Cells.Item(89, 11).Value = 42996
$ws.Cells.Item(89, 12).Value = 53657.5
$ws.Cells.Item(89, 13).Value = -37380
$ws.Cells.Item(89, 14).Value = -64889.5
$ws.Cells.Item(132, 8).Value = 1568.8939
$ws.Cells.Item(132, 9).Value = 1340.1803
$ws.Cells.Item(132, 11).Value = 4020.5409
$ws.Cells.Item(132, 13).Value = -1490.5409

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 2000
$ws.Cells.Item(39, 9).Value = 2000
$ws.Cells.Item(39, 11).Value = 6000
$ws.Cells.Item(39, 13).Value = -5706
$ws.Cells.Item(55, 8).Value = 135998
$ws.Cells.Item(55, 9).Value = 353333.34
$ws.Cells.Item(55, 10).Value = 5596.8
$ws.Cells.Item(55, 11).Value = 1060000.02
$ws.Cells.Item(55, 12).Value = 16790.4
$ws.Cells.Item(55, 13).Value = -1059823.02
$ws.Cells.Item(55, 14).Value = -17144.4
$ws.Cells.Item(80, 8).Value = 1281.6
$ws.Cells.Item(80, 9).Value = 1572
$ws.Cells.Item(80, 10).Value = 120
$ws.Cells.Item(80, 11).Value = 4716
$ws.Cells.Item(80, 12).Value = 360
$ws.Cells.Item(80, 13).Value = -3780
$ws.Cells.Item(80, 14).Value = -2232
$ws.Cells.Item(83, 8).Value = 1281.6
$ws.Cells.Item(83, 9).Value = 1572
$ws.Cells.Item(83, 10).Value = 120
$ws.Cells.Item(83, 11).Value = 14148
$ws.Cells.Item(83, 12).Value = 1080
$ws.Cells.Item(83, 13).Value = -9468
$ws.Cells.Item(83, 14).Value = -10440
$ws.Cells.Item(121, 8).Value = 849.4
$ws.Cells.Item(121, 9).Value = 800
$ws.Cells.Item(121, 10).Value = 861.75
$ws.Cells.Item(121, 11).Value = 2400
$ws.Cells.Item(121, 12).Value = 2585.25
$ws.Cells.Item(121, 13).Value = -1090
$ws.Cells.Item(121, 14).Value = -5205.25
$ws.Cells.Item(140, 8).Value = 1182.5834
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 1323.9166
$ws.Cells.Item(22, 10).Value = 1500
$ws.Cells.Item(22, 12).Value = 1500
$ws.Cells.Item(22, 14).Value = -2558
$ws.Cells.Item(80, 8).Value = 2280.2666
$ws.Cells.Item(80, 9).Value = 2009.4286
$ws.Cells.Item(80, 10).Value = 2517.25
$ws.Cells.Item(80, 11).Value = 2009.4286
$ws.Cells.Item(80, 12).Value = 2517.25
$ws.Cells.Item(80, 13).Value = -1011.4286
$ws.Cells.Item(80, 14).Value = -4513.25
$ws.Cells.Item(83, 8).Value = 2280.2666
$ws.Cells.Item(83, 9).Value = 2009.4286
$ws.Cells.Item(83, 10).Value = 2517.25
$ws.Cells.Item(83, 11).Value = 10047.143
$ws.Cells.Item(83, 12).Value = 12586.25
$ws.Cells.Item(83, 13).Value = -5055.143
$ws.Cells.Item(83, 14).Value = -22570.25
$ws.Cells.Item(113, 8).Value = 2999
$ws.Cells.Item(113, 9).Value = 2999
$ws.Cells.Item(113, 11).Value = 2999
$ws.Cells.Item(113, 13).Value = -829
$ws.Cells.Item(122, 8).Value = 102526.3
$ws.Cells.Item(122, 9).Value = 2883
$ws.Cells.Item(122, 10).Value = 501099.5
$ws.Cells.Item(122, 11).Value = 8649
$ws.Cells.Item(122, 12).Value = 1503298.5
$ws.Cells.Item(122, 13).Value = -6199
$ws.Cells.Item(122, 14).Value = -1508198.5
$ws.Cells.Item(123, 8).Value = 27462.666
$ws.Cells.Item(123, 10).Value = 27462.666
$ws.Cells.Item(123, 12).Value = 27462.666
$ws.Cells.Item(123, 14).Value = -32362.666
$ws.Cells.Item(126, 8).Value = 4506.5
$ws.Cells.Item(126, 9).Value = 2999
$ws.Cells.Item(126, 11).Value = 8997
$ws.Cells.Item(126, 13).Value = -6527

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3041.182
$ws.Cells.Item(46, 9).Value = 2100.5
$ws.Cells.Item(46, 10).Value = 4687.375
$ws.Cells.Item(46, 11).Value = 2100.5
$ws.Cells.Item(46, 12).Value = 4687.375
$ws.Cells.Item(46, 13).Value = -1912.5
$ws.Cells.Item(46, 14).Value = -5063.375
$ws.Cells.Item(55, 8).Value = 500
$ws.Cells.Item(55, 9).Value = 500
$ws.Cells.Item(55, 11).Value = 500
$ws.Cells.Item(55, 13).Value = -327
$ws.Cells.Item(60, 8).Value = 51580
$ws.Cells.Item(60, 9).Value = 42499.5
$ws.Cells.Item(60, 10).Value = 69741
$ws.Cells.Item(60, 11).Value = 42499.5
$ws.Cells.Item(60, 12).Value = 69741
$ws.Cells.Item(60, 13).Value = -41990.5
$ws.Cells.Item(60, 14).Value = -70759
$ws.Cells.Item(93, 8).Value = 1999.5
$ws.Cells.Item(93, 9).Value = 1999.5
$ws.Cells.Item(93, 11).Value = 1999.5
$ws.Cells.Item(93, 13).Value = -751.5
$ws.Cells.Item(122, 8).Value = 5517
$ws.Cells.Item(122, 9).Value = 6995
$ws.Cells.Item(122, 11).Value = 20985
$ws.Cells.Item(122, 13).Value = -18535
$ws.Cells.Item(136, 8).Value = 3999.75
$ws.Cells.Item(136, 9).Value = 3589.7778
$ws.Cells.Item(136, 10).Value = 5229.6665
$ws.Cells.Item(136, 11).Value = 10769.3334
$ws.Cells.Item(136, 12).Value = 15688.9995
$ws.Cells.Item(136, 13).Value = -8219.3334
$ws.Cells.Item(136, 14).Value = -20788.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 115362.11
$ws.Cells.Item(14, 9).Value = 127281.75
$ws.Cells.Item(14, 10).Value = 20005
$ws.Cells.Item(14, 11).Value = 127281.75
$ws.Cells.Item(14, 12).Value = 20005
$ws.Cells.Item(14, 13).Value = -127113.75
$ws.Cells.Item(14, 14).Value = -20341
$ws.Cells.Item(62, 8).Value = 8364.143
$ws.Cells.Item(62, 9).Value = 8049.5
$ws.Cells.Item(62, 10).Value = 8416.583000000001
$ws.Cells.Item(62, 11).Value = 8049.5
$ws.Cells.Item(62, 12).Value = 8416.583000000001
$ws.Cells.Item(62, 13).Value = -7425.5
$ws.Cells.Item(62, 14).Value = -9664.583000000001
$ws.Cells.Item(65, 8).Value = 8364.143
$ws.Cells.Item(65, 9).Value = 8049.5
$ws.Cells.Item(65, 10).Value = 8416.583000000001
$ws.Cells.Item(65, 11).Value = 40247.5
$ws.Cells.Item(65, 12).Value = 42082.915
$ws.Cells.Item(65, 13).Value = -37127.5
$ws.Cells.Item(65, 14).Value = -48322.915
$ws.Cells.Item(104, 8).Value = 21749.75
$ws.Cells.Item(104, 10).Value = 21749.75
$ws.Cells.Item(104, 12).Value = 21749.75
$ws.Cells.Item(104, 14).Value = -28737.75
$ws.Cells.Item(122, 8).Value = 3430.625
$ws.Cells.Item(122, 9).Value = 3430.625
$ws.Cells.Item(122, 11).Value = 10291.875
$ws.Cells.Item(122, 13).Value = -7841.875
$ws.Cells.Item(126, 8).Value = 2274.1875
$ws.Cells.Item(126, 9).Value = 1670.5
$ws.Cells.Item(126, 11).Value = 5011.5
$ws.Cells.Item(126, 13).Value = -2541.5

